# Apply the periodic cryptocurrency price/volume refresh described by the commit
# "Updated cryptos list on Tue Aug 29 06:51:47 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values in column D (e.g. "218.36") look like plain decimal numbers.
# The source data stores them as literal text (to preserve formats such as
# "26.113.98" or "0.0₅8160" elsewhere in the column), so force a Text number
# format on those specific cells before assigning, which keeps Range.Value from
# auto-converting them into Doubles.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D17", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D44", "D46", "D48", "D49", "D50", "D51", "D42")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Rows 2-41 and 44-51: refreshed Price (D) and/or Volume(1h) (E) readings
$ws.Range("D2").Value = "26.111.33"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.651.97"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "218.36"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.5208"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.2644"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "0.06339"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "20.36"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "0.07685"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "4.621"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "1.694.79"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "1.879.90"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "0.5597"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "0.0₅8157"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "65.36"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "26.110.04"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "4.629"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "10.50"
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").Value = "190.72"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "5.937"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "144.55"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "0.1191"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "7.218"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "15.93"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "1.510"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").Value = "0.05485"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").Value = "1.269"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "1.560"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").Value = "0.9490"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "2.784"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "2.398"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "0.5641"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "0.01577"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "5.862"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D44").Value = "100.90"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").Value = "1.791.99"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "57.72"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "0.0₈108"
$ws.Range("E47").Value = "  +5.51%  "
$ws.Range("D48").Value = "0.9999"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "0.4339"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "7.970"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "0.05174"
$ws.Range("E51").Value = "  -2.93%  "

# Rows 42 and 43 swapped ranking order (TrustWalletToken moved above Maker),
# each also carrying its own refreshed Price/Volume reading
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8329"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.029.78"
$ws.Range("E43").Value = "  -2.62%  "
